$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209, shifting existing rows 209:217 down to 210:218
$ws.Range("A209:T209").Insert()

# Copy formatting from the row below (the row that used to be 209, now 210) so the
# new row matches the style used throughout the data block (esp. date style on column D)
$ws.Range("A210:T210").Copy()
$ws.Range("A209:T209").PasteSpecial(-4122) # xlPasteFormats

# Populate the new row 209 with its data, following the same template as the
# surrounding rows (same market / region / product taxonomy), with the new
# record's own date, volume, prices and $/kg values.
$ws.Cells.Item(209, 1).Value = 4
$ws.Cells.Item(209, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(209, 3).Value = "Los Lagos"
$ws.Cells.Item(209, 4).Value = 44610
$ws.Cells.Item(209, 5).Value = 10
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100104
$ws.Cells.Item(209, 8).Value = "Frutos de pepita"
$ws.Cells.Item(209, 9).Value = 100104005
$ws.Cells.Item(209, 10).Value = "Pera"
$ws.Cells.Item(209, 11).Value = "Packham's Triumph"
$ws.Cells.Item(209, 12).Value = "Primera"
$ws.Cells.Item(209, 13).Value = 400
$ws.Cells.Item(209, 14).Value = 13000
$ws.Cells.Item(209, 15).Value = 13000
$ws.Cells.Item(209, 16).Value = 13000
$ws.Cells.Item(209, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(209, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(209, 19).Value = 867
$ws.Cells.Item(209, 20).Value = 15
